$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.842.70"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "2.923.66"
$ws.Range("E3").Value = "  +3.42%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'351.91"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "'111.71"
$ws.Range("E6").Value = "  -1.36%  "

$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "'39.25"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("E11").Value = "  +2.93%  "

$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("D13").Value = "'20.07"
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "3.392.80"
$ws.Range("E14").Value = "  +3.63%  "

$ws.Range("D15").Value = "'7.73"
$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").Value = "2.936.00"
$ws.Range("E16").Value = "  +3.70%  "

$ws.Range("D17").Value = "'0.979"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "51.938.80"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").Value = "'7.60"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("E20").Value = "  -4.04%  "

$ws.Range("D21").Value = "'14.17"
$ws.Range("E21").Value = "  +5.92%  "

$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").Value = "'71.11"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").Value = "'267.77"
$ws.Range("E24").Value = "  -0.54%  "

$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").Value = "'0.180"
$ws.Range("E26").Value = "  +11.09%  "

$ws.Range("D27").Value = "'26.92"
$ws.Range("E27").Value = "  +2.41%  "

$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").Value = "'7.25"
$ws.Range("E29").Value = "  +17.12%  "

$ws.Range("D30").Value = "'0.105"
$ws.Range("E30").Value = "  +16.60%  "

$ws.Range("D31").Value = "'10.57"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'37.13"
$ws.Range("E32").Value = "  -4.70%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'2.25"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").Value = "'6.21"
$ws.Range("E34").Value = "  +10.05%  "

$ws.Range("D35").Value = "'52.88"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("D36").Value = "'0.0452"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("D39").Value = "'18.62"
$ws.Range("E39").Value = "  -2.47%  "

$ws.Range("D40").Value = "'2.04"
$ws.Range("E40").Value = "  +1.49%  "

$ws.Range("E41").Value = "  +6.04%  "

$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("E43").Value = "  +5.66%  "

$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.50"
$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.173.69"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").Value = "'111.22"
$ws.Range("E48").Value = "  -8.73%  "

$ws.Range("E49").Value = "  +3.37%  "

$ws.Range("D50").Value = "'0.0348"
$ws.Range("E50").Value = "  +8.46%  "

$ws.Range("D51").Value = "'0.944"
$ws.Range("E51").Value = "  -5.63%  "

# Reset style for quote-prefixed numeric-looking text cells to avoid leftover quotePrefix style
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
